$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2024-04-12)
$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 9.826150383939911

# Row 3 (2024-04-07)
$ws.Range("B3").Value = 0.003994804209775715
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 4.443609181331619
